$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D-column price cells keep their text formatting (values contain
# multiple separators / trailing zeros that must not be parsed as numbers).
$dCells = @("D2","D3","D5","D6","D7","D8","D11","D12","D13","D14","D15","D16","D17","D18","D19","D20","D22","D23","D25","D26","D27","D29","D30","D31","D32","D33","D35","D36","D37","D38","D40","D41","D42","D43","D44","D46","D47","D48","D50","D51")
foreach ($ref in $dCells) { $ws.Range($ref).NumberFormat = "@" }

# Update price (D) values
$ws.Range("D2").Value = "70.436.33"
$ws.Range("D3").Value = "3.578.13"
$ws.Range("D5").Value = "586.66"
$ws.Range("D6").Value = "186.28"
$ws.Range("D7").Value = "3.564.47"
$ws.Range("D8").Value = "0.620"
$ws.Range("D11").Value = "0.650"
$ws.Range("D12").Value = "54.31"
$ws.Range("D13").Value = "0.0000313"
$ws.Range("D14").Value = "9.55"
$ws.Range("D15").Value = "4.141.96"
$ws.Range("D16").Value = "19.62"
$ws.Range("D17").Value = "70.427.23"
$ws.Range("D18").Value = "3.579.38"
$ws.Range("D19").Value = "12.43"
$ws.Range("D20").Value = "564.57"
$ws.Range("D22").Value = "1.02"
$ws.Range("D23").Value = "17.83"
$ws.Range("D25").Value = "4.92"
$ws.Range("D26").Value = "95.29"
$ws.Range("D27").Value = "11.60"
$ws.Range("D29").Value = "9.16"
$ws.Range("D30").Value = "32.20"
$ws.Range("D31").Value = "7.36"
$ws.Range("D32").Value = "12.50"
$ws.Range("D33").Value = "64.94"
$ws.Range("D35").Value = "3.37"
$ws.Range("D36").Value = "564.55"
$ws.Range("D37").Value = "0.422"
$ws.Range("D38").Value = "37.92"
$ws.Range("D40").Value = "0.0₃0775"
$ws.Range("D41").Value = "0.135"
$ws.Range("D42").Value = "3.362.96"
$ws.Range("D43").Value = "3.38"
$ws.Range("D44").Value = "3.05"
$ws.Range("D46").Value = "0.0445"
$ws.Range("D47").Value = "2.96"
$ws.Range("D48").Value = "9.43"
$ws.Range("D50").Value = "0.999"
$ws.Range("D51").Value = "1.47"

# Update other values (volume %, coin name, link)
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("E3").Value = "  -0.07%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("E5").Value = "  +1.48%  "
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("E7").Value = "  -0.29%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("E10").Value = "  +11.35%  "
$ws.Range("E11").Value = "  -0.32%  "
$ws.Range("E12").Value = "  -1.53%  "
$ws.Range("E13").Value = "  +2.17%  "
$ws.Range("E14").Value = "  +0.05%  "
$ws.Range("E15").Value = "  -0.28%  "
$ws.Range("E16").Value = "  -0.35%  "
$ws.Range("E17").Value = "  +0.56%  "
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("E19").Value = "  -1.58%  "
$ws.Range("E20").Value = "  +15.28%  "
$ws.Range("E21").Value = "  -0.32%  "
$ws.Range("E22").Value = "  -1.63%  "
$ws.Range("E23").Value = "  -6.67%  "
$ws.Range("E24").Value = "  +5.62%  "
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("E27").Value = "  -2.79%  "
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("E29").Value = "  -1.82%  "
$ws.Range("E30").Value = "  +1.80%  "
$ws.Range("E31").Value = "  -4.98%  "
$ws.Range("E32").Value = "  +3.43%  "
$ws.Range("E33").Value = "  -3.05%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("E35").Value = "  +2.88%  "
$ws.Range("E36").Value = "  -1.60%  "
$ws.Range("E37").Value = "  +0.27%  "
$ws.Range("E38").Value = "  -2.00%  "
$ws.Range("E39").Value = "  +0.25%  "
$ws.Range("E40").Value = "  -2.34%  "
$ws.Range("E41").Value = "  -0.56%  "
$ws.Range("E42").Value = "  +4.39%  "
$ws.Range("E43").Value = "  -2.21%  "
$ws.Range("E44").Value = "  -3.39%  "
$ws.Range("E45").Value = "  +2.32%  "
$ws.Range("E46").Value = "  +0.16%  "
$ws.Range("E47").Value = "  -3.69%  "
$ws.Range("E48").Value = "  -1.69%  "
$ws.Range("E49").Value = "  +0.41%  "
$ws.Range("B50").Value = "FirstDigitalUSD"
$ws.Range("C50").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("E50").Value = "  -0.04%  "
$ws.Range("B51").Value = "OceanProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/aAKLSV5-0+oceanprotocol-ocean"
$ws.Range("E51").Value = "  -9.61%  "
